$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (participant 103): PhotosUploaded? changed from "n" to "y"
$ws.Range("E4").Value = "y"

# New row 5 for participant 105
$ws.Range("A4").Copy()
$ws.Range("A5").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A5").Value = (Get-Date -Year 2025 -Month 10 -Day 22 -Hour 0 -Minute 0 -Second 0)
$ws.Range("B5").Value = 105
$ws.Range("C5").Value = "R_34L7XTVaE312VgJ"
$ws.Range("D5").Value = "R_6pffiNtHXao5FYL"
$ws.Range("E5").Value = "n"

# Update selection to match diff (active cell D7)
$ws.Range("D7").Select()
